$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.721.61'
$ws.Range("E2").Value = '  +5.72%  '

$ws.Range("D3").Value = '2.261.35'
$ws.Range("E3").Value = '  +4.47%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.07'
$ws.Range("E5").Value = '  +3.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.643'
$ws.Range("E6").Value = '  +3.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.07'
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +4.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.85'
$ws.Range("E10").Value = '  +3.03%  '

$ws.Range("E11").Value = '  +5.53%  '

$ws.Range("E12").Value = '  +2.17%  '

$ws.Range("D13").Value = '2.596.29'
$ws.Range("E13").Value = '  +4.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.14'
$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.89'
$ws.Range("E15").Value = '  +3.71%  '

$ws.Range("E16").Value = '  +1.77%  '

$ws.Range("E17").Value = '  +3.34%  '

$ws.Range("D18").Value = '2.288.46'
$ws.Range("E18").Value = '  +5.40%  '

$ws.Range("D19").Value = '41.625.41'

$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("E20").Value = '  +10.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.87'
$ws.Range("E21").Value = '  +4.19%  '

$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.23'
$ws.Range("E23").Value = '  +9.89%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  +3.65%  '

$ws.Range("E26").Value = '  +0.54%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.87'
$ws.Range("E27").Value = '  +3.66%  '

$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.149'
$ws.Range("E28").Value = '  +5.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.65'
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.54'
$ws.Range("E30").Value = '  +3.40%  '

$ws.Range("E31").Value = '  +2.23%  '

$ws.Range("E32").Value = '  +5.45%  '

$ws.Range("E33").Value = '  +3.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.12'
$ws.Range("E34").Value = '  +8.37%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.83'
$ws.Range("E35").Value = '  +4.37%  '

$ws.Range("E36").Value = '  +2.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.91'
$ws.Range("E37").Value = '  -2.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.87'
$ws.Range("E38").Value = '  +7.25%  '

$ws.Range("E39").Value = '  +1.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000268'
$ws.Range("E40").Value = '  +65.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.09'
$ws.Range("E41").Value = '  +19.55%  '

$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("E43").Value = '  +5.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.73'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.88'
$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0989'
$ws.Range("E46").Value = '  +6.75%  '

$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.24'
$ws.Range("E47").Value = '  +1.94%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '17.64'
$ws.Range("E48").Value = '  -0.99%  '

$ws.Range("D49").Value = '1.511.54'
$ws.Range("E49").Value = '  -0.64%  '

$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("E51").Value = '  -0.86%  '
